# Applies hybrid bold + color (2C3E50) highlighting to quantitative
# metrics (percentages, dollar amounts, large numbers) inside specific
# bullet/summary paragraphs, matching the target diff.
#
# For each target paragraph, the plain single run is split into
# alternating plain / bold+colored runs around each metric token by
# scoping Word's Find.Execute to the paragraph (and to the remaining
# tail of it, so repeated tokens such as "87%" are matched left to
# right and only once per occurrence requested).

$d = $word.ActiveDocument

# Word stores RGB colors in OLE COLORREF (BGR-packed) integer form.
# 0x2C3E50 (R=2C,G=3E,B=50) -> BGR packed = 0x503E2C = 5258796
$metricColor = 5258796

function Apply-MetricHighlights($Paragraph, $ExpectedSnippet, $Metrics) {
    $pRange = $Paragraph.Range.Duplicate
    $paraEnd = $pRange.End
    $cursor = $pRange.Start

    if ($pRange.Text.IndexOf($ExpectedSnippet) -lt 0) {
        Write-Output "WARNING: paragraph did not contain expected snippet: $ExpectedSnippet"
    }

    foreach ($metric in $Metrics) {
        $searchRange = $d.Range($cursor, $paraEnd)
        $found = $searchRange.Find.Execute(
            $metric, $false, $false, $false, $false, $false,
            $true, 1, $false, "", 0)

        if ($found) {
            $searchRange.Font.Bold = $true
            $searchRange.Font.Color = $metricColor
            $cursor = $searchRange.End
        }
        else {
            Write-Output "WARNING: metric not found: $metric"
        }
    }
}

# 1) "...demographic classification accuracy from 23% to 64%"
Apply-MetricHighlights $d.Paragraphs.Item(9) "23% to 64%" @("23%", "64%")

# 2) "Achieved 87% ... industry standard of 71%, ... from \xB14.2% to \xB12.1%"
Apply-MetricHighlights $d.Paragraphs.Item(11) "87% prediction accuracy" @(
    "87%", "71%", ([char]0x00B1 + "4.2%"), ([char]0x00B1 + "2.1%"))

# 3) "...bids from 1,200 vendors..."
Apply-MetricHighlights $d.Paragraphs.Item(31) "1,200 vendors" @("1,200")

# 4) "...became the $400M Polling Consortium Database ... valued at $1B+"
Apply-MetricHighlights $d.Paragraphs.Item(46) "$400M Polling Consortium" @("$400M", "$1B")

# 5) "Algorithm reduced mapping costs by 73.5%, saving ... $4.7M"
Apply-MetricHighlights $d.Paragraphs.Item(63) "Algorithm reduced mapping costs" @("73.5%", "$4.7M")

# 6) "Achieved 87% prediction accuracy ... industry standard of 71%"
Apply-MetricHighlights $d.Paragraphs.Item(65) "87% prediction accuracy" @("87%", "71%")
